# Updated cryptos list on Sat Sep 16 15:33:49 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as TEXT (matches source data which is all inline
# strings, even for number-looking prices) without leaving the cells
# number format permanently changed - stash/restore the original Style.
function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue "D2" "26.597.64"
Set-TextValue "E2" "  +0.82%  "

Set-TextValue "D3" "1.640.14"
Set-TextValue "E3" "  +1.11%  "

Set-TextValue "E4" "  -0.11%  "

Set-TextValue "D5" "214.48"
Set-TextValue "E5" "  +1.02%  "

Set-TextValue "E6" "  +1.56%  "

Set-TextValue "E7" "  -0.11%  "

Set-TextValue "E8" "  +1.27%  "

Set-TextValue "E9" "  +0.70%  "

Set-TextValue "D10" "19.07"
Set-TextValue "E10" "  +0.95%  "

Set-TextValue "D11" "0.0841"
Set-TextValue "E11" "  +0.05%  "

Set-TextValue "E12" "  +1.04%  "

Set-TextValue "D13" "1.633.36"
Set-TextValue "E13" "  +0.23%  "

Set-TextValue "E14" "  +1.87%  "

Set-TextValue "E15" "  +1.59%  "

Set-TextValue "E16" "  +1.20%  "

Set-TextValue "D17" "26.615.05"
Set-TextValue "E17" "  +0.87%  "

Set-TextValue "E18" "  +0.46%  "

Set-TextValue "D19" "215.14"
Set-TextValue "E19" "  +0.86%  "

Set-TextValue "E21" "  +0.97%  "

Set-TextValue "D22" "6.24"
Set-TextValue "E22" "  +0.69%  "

Set-TextValue "D23" "9.44"
Set-TextValue "E23" "  +1.84%  "

Set-TextValue "E24" "  +12.92%  "

Set-TextValue "D25" "145.00"
Set-TextValue "E25" "  -1.76%  "

Set-TextValue "E26" "  -0.08%  "

Set-TextValue "E27" "  +0.10%  "

Set-TextValue "E28" "  +4.39%  "

Set-TextValue "D29" "15.66"
Set-TextValue "E29" "  +0.92%  "

Set-TextValue "E30" "  +1.22%  "

Set-TextValue "E31" "  +1.48%  "

Set-TextValue "D32" "3.36"
Set-TextValue "E32" "  +1.38%  "

Set-TextValue "E33" "  +1.66%  "

Set-TextValue "D34" "1.278.72"
Set-TextValue "E34" "  +5.67%  "

Set-TextValue "E35" "  +2.83%  "

Set-TextValue "D36" "2.41"
Set-TextValue "E36" "  +1.16%  "

Set-TextValue "E37" "  +2.97%  "

Set-TextValue "D38" "0.531"
Set-TextValue "E38" "  +6.42%  "

Set-TextValue "E39" "  +3.59%  "

Set-TextValue "E40" "  -0.08%  "

Set-TextValue "E41" "  +2.22%  "

Set-TextValue "E42" "  -0.12%  "

Set-TextValue "E43" "  +0.78%  "

Set-TextValue "D44" "1.778.68"
Set-TextValue "E44" "  +1.07%  "

Set-TextValue "D45" "91.32"
Set-TextValue "E45" "  -1.30%  "

Set-TextValue "D46" "59.29"
Set-TextValue "E46" "  +8.71%  "

Set-TextValue "E47" "  +1.14%  "

Set-TextValue "B48" "Cronos"
Set-TextValue "C48" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D48" "0.0514"
Set-TextValue "E48" "  +0.91%  "

Set-TextValue "B49" "EnergySwap"
Set-TextValue "C49" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D49" "7.72"
Set-TextValue "E49" "  +1.48%  "

Set-TextValue "B50" "Algorand"
Set-TextValue "C50" "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D50" "0.0961"
Set-TextValue "E50" "  +1.28%  "

Set-TextValue "B51" "Mantle"
Set-TextValue "C51" "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D51" "0.405"
Set-TextValue "E51" "  -0.59%  "
